# Generate Report for handoff
# Adds a new "ready for handoff" row for file
# c6ff638b-4995-4d5e-92c6-05d7e422a72c.md to the Overview, zh-cn and de-de
# sheets, pushing the existing ".localization-config" row down by one.

$wb = $excel.ActiveWorkbook

$newFile = "c6ff638b-4995-4d5e-92c6-05d7e422a72c.md"
$newHash = "9f854d74d4e441bcc99ec294df790fca43910e5b"
$zhXlfName = "c6ff638b-4995-4d5e-92c6-05d7e422a72c.$newHash.zh-cn.xlf"
$deXlfName = "c6ff638b-4995-4d5e-92c6-05d7e422a72c.$newHash.de-de.xlf"
$zhHandoffDateTime = "2016-01-20 07:33:20"
$deHandoffDateTime = "2016-01-20 07:33:31"
$epoch = "0001-01-01 00:00:00"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/a7962610fbdb89c950679a70a424b3f6962e7d16/e2e/$newFile"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

# Capture the existing hyperlinks (in sheet order) before we touch the grid,
# so that we can rebuild them afterwards with correct ranges / rIds, in the
# correct final order.
$oldLinks = @()
foreach ($h in $wsO.Hyperlinks) {
    $oldLinks += , @($h.Range.Row, $h.Range.Column, $h.Address, $h.TextToDisplay)
}

# Insert a new row above the ".localization-config" row (row 8), copying the
# formatting from the row above it so the new row matches the table style.
$wsO.Rows.Item(7).Copy()
$wsO.Rows.Item(8).Insert()

$wsO.Range("A8").Value = $newFile
$wsO.Range("B8").Value = "Ready for handoff"
$wsO.Range("C8").Value = "Ready for handoff"

# Rebuild all hyperlinks in final document order: links before row 8 stay
# put, then the new row-8 link, then the old row-8 link (now row 9).
$wsO.Hyperlinks.Delete()
foreach ($link in $oldLinks) {
    $row = $link[0]
    $col = $link[1]
    if ($row -eq 8) {
        $wsO.Hyperlinks.Add($wsO.Range("A8"), $mdUrl, "", "", $newFile) | Out-Null
        $row = 9
    }
    $cell = $wsO.Cells.Item($row, $col)
    $wsO.Hyperlinks.Add($cell, $link[2], "", "", $link[3]) | Out-Null
}

# ---------------------------------------------------------------------------
# Helper that performs the equivalent edit on a per-language detail sheet
# (zh-cn / de-de), which share the same layout.
# ---------------------------------------------------------------------------
function Update-LangSheet($sheetName, $xlfName, $xlfUrl, $handoffDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $oldLinks = @()
    foreach ($h in $ws.Hyperlinks) {
        $oldLinks += , @($h.Range.Row, $h.Range.Column, $h.Address, $h.TextToDisplay)
    }

    $ws.Rows.Item(7).Copy()
    $ws.Rows.Item(8).Insert()

    $ws.Range("A8").Value = $newFile
    $ws.Range("B8").Value = "Ready for handoff"
    $ws.Range("C8").Value = $xlfName
    $ws.Range("D8").Value = $handoffDateTime
    $ws.Range("E8").ClearContents()
    $ws.Range("F8").ClearContents()
    $ws.Range("G8").Value = $epoch
    $ws.Range("H8").Value = "Include"

    $ws.Range("D9").Value = $epoch
    $ws.Range("G9").Value = $epoch
    $ws.Range("H9").Value = "Ignored"

    $ws.Hyperlinks.Delete()
    foreach ($link in $oldLinks) {
        $row = $link[0]
        $col = $link[1]
        if ($row -eq 8) {
            $ws.Hyperlinks.Add($ws.Range("A8"), $mdUrl, "", "", $newFile) | Out-Null
            $ws.Hyperlinks.Add($ws.Range("C8"), $xlfUrl, "", "", $xlfName) | Out-Null
            $row = 9
        }
        $cell = $ws.Cells.Item($row, $col)
        $ws.Hyperlinks.Add($cell, $link[2], "", "", $link[3]) | Out-Null
    }
}

Update-LangSheet "zh-cn" $zhXlfName $zhXlfUrl $zhHandoffDateTime
Update-LangSheet "de-de" $deXlfName $deXlfUrl $deHandoffDateTime
